# Generate Report for Handback
# Updates the localization-status workbook: marks the "a.md" entries as
# handed back (in sync with en-US), records the handback xlf file + datetime
# for both the zh-cn and de-de target languages, and links the new
# "Latest Target File" cell back to the source markdown file on GitHub.

$wb = $excel.ActiveWorkbook

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b97c3be6944add41f1b94f0a16c52cc898b2e5c9/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: both zh-cn and de-de columns for a.md flip from
# "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the Status columns so the longer text fits (matches generated report).
$wsOverview.Range("E1").ColumnWidth = 29
$wsOverview.Range("F1").ColumnWidth = 29

# ---------------------------------------------------------------------
# zh-cn sheet: row 2 (a.md) and row 3 (b.md) both get a Latest Target
# File hyperlink back to a.md, a Latest Handback File name, and a Latest
# Handback DateTime. Status also flips to "Handed back" (shares the same
# string as the Overview sheet's status cells).
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aMdUrl, "", "", "a.md")
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276

$wsZh.Range("I3").Value = "a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aMdUrl, "", "", "a.md")
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-28 02:37:11"
$wsZh.Range("K3").Value = "2016-08-28 02:37:11"

$wsZh.Range("C1").ColumnWidth = 29
$wsZh.Range("J1").ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, but with the de-de handback
# xlf file name and its own handback timestamp.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aMdUrl, "", "", "a.md")
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276

$wsDe.Range("I3").Value = "a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aMdUrl, "", "", "a.md")
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-28 02:37:19"
$wsDe.Range("K3").Value = "2016-08-28 02:37:19"

$wsDe.Range("C1").ColumnWidth = 29
$wsDe.Range("J1").ColumnWidth = 39.1666666666667

Write-Output "Handback report generated"
